$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# New "waitSec" column (G) on the second sheet.
$ws2.Range("G1").Value = "waitSec"
$ws2.Range("G2").Value = 5

# Make sheet 2 the active tab / selected sheet, with H12 selected.
$ws2.Activate()
$ws2.Range("H12").Select()
